$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13: fill in E13/F13, recompute G13 (shared formula already present),
#     add H13 and I13 (new shared formula sibling of I11:I12 group) ---
$ws.Range("E13").Value = 417
$ws.Range("F13").Value = 448
$ws.Range("H13").Value = 4
$ws.Range("I13").Formula = "=G13/H13"

# --- K1: extend AVERAGE range to include row 13 ---
$ws.Range("K1").Formula = "=AVERAGE(G2:G13)"

# --- K2: reference F13 instead of F12 ---
$ws.Range("K2").Formula = "=(B15-F13)/K1"

# --- New label + median formula on row 4 ---
$ws.Range("J4").Value = "Медиана стр/пом"
$ws.Range("K4").Formula = "=MEDIAN(I2:I13)"
$ws.Range("K4").NumberFormat = $ws.Range("K1").NumberFormat

# --- Update sheet view selection/scroll position ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("K5").Select()

$wb.Save()
